$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full data grid (header row + 5 data rows, columns A-W)
$arr = New-Object 'object[,]' 6,23
$arr[0,0] = "productIds"
$arr[0,1] = "MSE_no_transfer"
$arr[0,2] = "MSE_transfer_basic"
$arr[0,3] = "MSE_transfer_coral"
$arr[0,4] = "MSE_transfer_sa"
$arr[0,5] = "MSE_transfer_bw"
$arr[0,6] = "MSE_transfer_nnw"
$arr[0,7] = "MAE_no_transfer"
$arr[0,8] = "MAE_transfer_basic"
$arr[0,9] = "MAE_transfer_coral"
$arr[0,10] = "MAE_transfer_sa"
$arr[0,11] = "MAE_transfer_bw"
$arr[0,12] = "MAE_transfer_nnw"
$arr[0,13] = "MSE_diff_basic"
$arr[0,14] = "MSE_transfer_coral"
$arr[0,15] = "MSE_diff_sa"
$arr[0,16] = "MSE_diff_bw"
$arr[0,17] = "MSE_diff_nnw"
$arr[0,18] = "MAE_diff_basic"
$arr[0,19] = "MAE_transfer_coral"
$arr[0,20] = "MAE_diff_sa"
$arr[0,21] = "MAE_diff_bw"
$arr[0,22] = "MAE_diff_nnw"
$arr[1,0] = "101-120"
$arr[1,1] = 1.300493167645185
$arr[1,2] = 1.000568226643996
$arr[1,3] = 1.367077284212216
$arr[1,4] = 1.300356127336072
$arr[1,5] = 1.053588503088057
$arr[1,6] = 1.13702385724122
$arr[1,7] = 0.6937269902233257
$arr[1,8] = 0.6807548401453892
$arr[1,9] = 0.9693148174296417
$arr[1,10] = 0.906898140459924
$arr[1,11] = 0.7776672221506405
$arr[1,12] = 0.5497553639747121
$arr[1,13] = -0.299924941001189
$arr[1,14] = 0.06658411656703089
$arr[1,15] = -0.0001370403091132122
$arr[1,16] = -0.2469046645571282
$arr[1,17] = -0.1634693104039648
$arr[1,18] = -0.01297215007793651
$arr[1,19] = 0.2755878272063159
$arr[1,20] = 0.2131711502365983
$arr[1,21] = 0.08394023192731481
$arr[1,22] = -0.1439716262486136
$arr[2,0] = "121-140"
$arr[2,1] = 0.7858089996152173
$arr[2,2] = 0.9959746491911219
$arr[2,3] = 2.610539538751953
$arr[2,4] = 3.541888404581724
$arr[2,5] = 2.515436924917919
$arr[2,6] = 2.855604247748933
$arr[2,7] = 0.5457719070059673
$arr[2,8] = 0.664341957982579
$arr[2,9] = 1.129524704241595
$arr[2,10] = 1.2379582477458
$arr[2,11] = 1.06147489521414
$arr[2,12] = 0.9761645511911143
$arr[2,13] = 0.2101656495759046
$arr[2,14] = 1.824730539136736
$arr[2,15] = 2.756079404966506
$arr[2,16] = 1.729627925302702
$arr[2,17] = 2.069795248133715
$arr[2,18] = 0.1185700509766117
$arr[2,19] = 0.5837527972356272
$arr[2,20] = 0.6921863407398331
$arr[2,21] = 0.5157029882081726
$arr[2,22] = 0.430392644185147
$arr[3,0] = "141-160"
$arr[3,1] = 1.996350390516476
$arr[3,2] = 1.03790889963132
$arr[3,3] = 5.202889446466702
$arr[3,4] = 5.023626319614601
$arr[3,5] = 5.207028709953064
$arr[3,6] = 6.186496171442877
$arr[3,7] = 0.8411752455834071
$arr[3,8] = 0.6918812668648573
$arr[3,9] = 1.419065534568517
$arr[3,10] = 1.409730946063732
$arr[3,11] = 1.38916137543472
$arr[3,12] = 1.411307742520673
$arr[3,13] = -0.9584414908851564
$arr[3,14] = 3.206539055950226
$arr[3,15] = 3.027275929098125
$arr[3,16] = 3.210678319436588
$arr[3,17] = 4.1901457809264
$arr[3,18] = -0.1492939787185498
$arr[3,19] = 0.5778902889851102
$arr[3,20] = 0.5685557004803249
$arr[3,21] = 0.5479861298513128
$arr[3,22] = 0.5701324969372662
$arr[4,0] = "161-180"
$arr[4,1] = 0.9698881290922886
$arr[4,2] = 1.032012595198869
$arr[4,3] = 13.67423265531376
$arr[4,4] = 13.96535841057239
$arr[4,5] = 12.51438321513197
$arr[4,6] = 13.09847841109867
$arr[4,7] = 0.3799394289281765
$arr[4,8] = 0.4221985870863751
$arr[4,9] = 1.488683363065025
$arr[4,10] = 1.599282058080037
$arr[4,11] = 1.357240006472364
$arr[4,12] = 1.267557036725364
$arr[4,13] = 0.06212446610658018
$arr[4,14] = 12.70434452622147
$arr[4,15] = 12.99547028148011
$arr[4,16] = 11.54449508603968
$arr[4,17] = 12.12859028200638
$arr[4,18] = 0.04225915815819864
$arr[4,19] = 1.108743934136848
$arr[4,20] = 1.219342629151861
$arr[4,21] = 0.9773005775441874
$arr[4,22] = 0.8876176077971878
$arr[5,0] = "181-200"
$arr[5,1] = 1.51415165254623
$arr[5,2] = 1.011790435927412
$arr[5,3] = 7.869639472631155
$arr[5,4] = 7.510143915203348
$arr[5,5] = 7.900601669343906
$arr[5,6] = 8.277631412360071
$arr[5,7] = 0.5200546449324573
$arr[5,8] = 0.4728316765745479
$arr[5,9] = 1.304881907658572
$arr[5,10] = 1.367534295413275
$arr[5,11] = 1.243400161246943
$arr[5,12] = 1.156648651897284
$arr[5,13] = -0.5023612166188185
$arr[5,14] = 6.355487820084925
$arr[5,15] = 5.995992262657118
$arr[5,16] = 6.386450016797675
$arr[5,17] = 6.763479759813841
$arr[5,18] = -0.0472229683579094
$arr[5,19] = 0.7848272627261142
$arr[5,20] = 0.8474796504808176
$arr[5,21] = 0.7233455163144858
$arr[5,22] = 0.6365940069648269

$ws.Range("A1:W6").Value = $arr

# Apply the header style (bold, bordered, centered) to the newly added header cells L1:W1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("L1:W1").PasteSpecial(-4122) | Out-Null

Write-Output "Edit applied successfully"
